$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    'Bitcoin',
    'Ethereum',
    'BNB',
    'XRP',
    'Solana',
    'Cardano',
    'Dogecoin',
    'TRON',
    'Polygon',
    'Chainlink',
    'Toncoin',
    'Avalanche',
    'Polkadot',
    'Wrapped Bitcoin',
    'Litecoin',
    'Shiba Inu',
    'Bitcoin Cash',
    'UNUS SED LEO',
    'Cosmos',
    'OKB',
    'Stellar',
    'Monero',
    'Ethereum Classic',
    'Cronos',
    'Kaspa'
)

$symbols = @(
    'BTC-USD',
    'ETH-USD',
    'BNB-USD',
    'XRP-USD',
    'SOL-USD',
    'ADA-USD',
    'DOGE-USD',
    'TRX-USD',
    'MATIC-USD',
    'LINK-USD',
    'TON-USD',
    'AVAX-USD',
    'DOT-USD',
    'WBTC-USD',
    'LTC-USD',
    'SHIB-USD',
    'BCH-USD',
    'LEO-USD',
    'ATOM-USD',
    'OKB-USD',
    'XLM-USD',
    'XMR-USD',
    'ETC-USD',
    'CRO-USD',
    'KAS-USD'
)

$caps = @(
    708479008611.6012,
    242091520445.6097,
    38255025955.62946,
    34204919250.52261,
    26411865956.14902,
    13006021840.74045,
    10523736344.92889,
    9236740082.654627,
    8378001144.125881,
    8218904526.219995,
    7867658202.653782,
    7018185773.453343,
    6752494904.928741,
    5922102614.900625,
    5318515548.766405,
    5059433485.511156,
    4556617892.441616,
    3707281230.704166,
    3575037140.071479,
    3561330600.496291,
    3356540063.866071,
    3000328346.07416,
    2795690272.088627,
    2563723725.813462,
    2514361673.708901
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $symbols[$i]
    $ws.Cells.Item($row, 3).Value = $caps[$i]
}

Write-Output "Done writing data rows."